$d = $word.ActiveDocument

$pairs = @(
    @("4+19=23", "17+35=52"),
    @("30-9=21", "16-9=7"),
    @("70-69=1", "38+53=91"),
    @("46+19=65", "56-8=48"),
    @("43-17=26", "81-79=2"),
    @("46-8=38", "25+27=52"),
    @("17+15=32", "69+15=84"),
    @("19+65=84", "6+39=45"),
    @("17+75=92", "42-19=23"),
    @("81-33=48", "90-5=85"),
    @("61-24=37", "68+6=74"),
    @("54+17=71", "8+9=17"),
    @("90-29=61", "8+59=67"),
    @("39+6=45", "34+38=72"),
    @("26+6=32", "83-18=65"),
    @("72-63=9", "55-37=18"),
    @("8+16=24", "32+19=51"),
    @("35+36=71", "57+5=62"),
    @("74+19=93", "43-14=29"),
    @("5+18=23", "64-5=59"),
    @("77+18=95", "56-28=28"),
    @("28+19=47", "6+76=82"),
    @("23+38=61", "28+36=64"),
    @("87-29=58", "7+24=31"),
    @("45-17=28", "16+59=75"),
    @("8+8=16", "43+48=91"),
    @("37-18=19", "74-49=25"),
    @("87-28=59", "20-18=2"),
    @("7+19=26", "39+28=67"),
    @("34-28=6", "32-18=14"),
    @("56-49=7", "85-18=67"),
    @("83+9=92", "46-17=29"),
    @("41-33=8", "95-46=49"),
    @("44-29=15", "52-23=29"),
    @("56+37=93", "80-77=3"),
    @("47-28=19", "37-8=29"),
    @("73-66=7", "63-28=35"),
    @("38+57=95", "2+79=81"),
    @("6+88=94", "55+18=73"),
    @("84-39=45", "27+56=83"),
    @("72-25=47", "78+6=84"),
    @("93-29=64", "18+36=54"),
    @("17+25=42", "7+34=41"),
    @("49+12=61", "30-12=18"),
    @("29+47=76", "46+8=54"),
    @("56+26=82", "60-58=2"),
    @("6+77=83", "7+35=42"),
    @("8+83=91", "76-39=37"),
    @("71-45=26", "17+49=66"),
    @("97-39=58", "26-8=18"),
    @("40-23=17", "70-47=23"),
    @("76-27=49", "23+59=82"),
    @("95-78=17", "54+39=93"),
    @("9+45=54", "56+8=64"),
    @("19+75=94", "38+55=93"),
    @("60-39=21", "91-34=57"),
    @("19+16=35", "83-68=15"),
    @("97-49=48", "90-46=44"),
    @("34+7=41", "71-46=25"),
    @("22+49=71", "89+9=98"),
    @("92-67=25", "34+28=62"),
    @("19+42=61", "17+14=31"),
    @("75-59=16", "28+53=81"),
    @("42-16=26", "21-6=15"),
    @("48-9=39", "17+54=71"),
    @("19+12=31", "54+19=73"),
    @("49+3=52", "51-3=48"),
    @("18+23=41", "5+7=12"),
    @("51-17=34", "84-5=79"),
    @("73-26=47", "12-4=8"),
    @("9+16=25", "85-78=7"),
    @("53-39=14", "68-39=29"),
    @("63-18=45", "80-16=64"),
    @("47+39=86", "30-5=25"),
    @("8+18=26", "38+5=43"),
    @("45+16=61", "8+14=22"),
    @("70-49=21", "75-28=47"),
    @("26+46=72", "61-57=4"),
    @("56+16=72", "7+74=81"),
    @("46-29=17", "78+13=91"),
    @("60-19=41", "82-14=68"),
    @("18+63=81", "35-27=8"),
    @("55+8=63", "31-6=25"),
    @("43-16=27", "55-19=36"),
    @("36+47=83", "26+28=54"),
    @("93-84=9", "18+3=21"),
    @("33-15=18", "63+18=81"),
    @("74-19=55", "17+68=85"),
    @("19+25=44", "58+34=92"),
    @("50-37=13", "21-17=4"),
    @("29+8=37", "64-18=46"),
    @("83-56=27", "33+18=51"),
    @("34+59=93", "83-29=54"),
    @("6+37=43", "6+49=55"),
    @("6+26=32", "15+56=71"),
    @("38+25=63", "61-6=55"),
    @("93-55=38", "59+17=76"),
    @("50-5=45", "62-19=43"),
    @("61-56=5", "47+38=85"),
    @("6+16=22", "8+36=44"),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done replacing pairs"